# Update "Pertanggal 3 Februari 2026 19:05 WIB"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Supplier row (row 9): E9 and F9 change from 1423 to 1520.
# Dependent formulas (G9, E14, F14, G14) will recalculate automatically.
$ws.Range("E9").Value = 1520
$ws.Range("F9").Value = 1520

# Hide row 12 (the "Project" row)
$ws.Rows(12).Hidden = $true

# Update the active selection to F10
$ws.Range("F10").Select()
